$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.050088250651209
$ws.Range("D2").Value = 1.053871971731013
$ws.Range("E2").Value = 1.04720179193776
$ws.Range("F2").Value = 1.062770925183947
$ws.Range("I2").Value = 1.044658817302352
$ws.Range("J2").Value = 1.055123255556531
$ws.Range("K2").Value = 1.056616743001365
$ws.Range("L2").Value = 1.049965090369585
$ws.Range("M2").Value = 1.065491370679172
$ws.Range("N2").Value = 1.05662165160375
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.05124482394254
$ws.Range("D3").Value = 1.054789724880777
$ws.Range("E3").Value = 1.048193626022795
$ws.Range("F3").Value = 1.063921443156906
$ws.Range("I3").Value = 1.044993886281747
$ws.Range("J3").Value = 1.055928358228258
$ws.Range("K3").Value = 1.057347359245543
$ws.Range("L3").Value = 1.050768253944747
$ws.Range("M3").Value = 1.066455926211547
$ws.Range("N3").Value = 1.057427897613618
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.051993137579391
$ws.Range("D4").Value = 1.055383508792437
$ws.Range("E4").Value = 1.048835634532259
$ws.Range("F4").Value = 1.064666326931831
$ws.Range("I4").Value = 1.045209467676163
$ws.Range("J4").Value = 1.05644869541789
$ws.Range("K4").Value = 1.057819423044739
$ws.Range("L4").Value = 1.051287566796844
$ws.Range("M4").Value = 1.067079895430331
$ws.Range("N4").Value = 1.057948973741743
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.052307713814832
$ws.Range("D5").Value = 1.055633120500823
$ws.Range("E5").Value = 1.049105589201874
$ws.Range("F5").Value = 1.064979577595884
$ws.Range("I5").Value = 1.045299803654903
$ws.Range("J5").Value = 1.056667297675913
$ws.Range("K5").Value = 1.058017712913057
$ws.Range("L5").Value = 1.051505793462157
$ws.Range("M5").Value = 1.067342173618817
$ws.Range("N5").Value = 1.058167886440046
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.052360531713034
$ws.Range("D6").Value = 1.055675030528768
$ws.Range("E6").Value = 1.04915091897119
$ws.Range("F6").Value = 1.065032179730735
$ws.Range("I6").Value = 1.045314954202376
$ws.Range("J6").Value = 1.056703993321411
$ws.Range("K6").Value = 1.05805099694443
$ws.Range("L6").Value = 1.051542429268093
$ws.Range("M6").Value = 1.067386209038826
$ws.Range("N6").Value = 1.05820463419757
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.051997341019242
$ws.Range("D7").Value = 1.055386844175164
$ws.Range("E7").Value = 1.04883924146606
$ws.Range("F7").Value = 1.064670512202897
$ws.Range("I7").Value = 1.045210675906652
$ws.Range("J7").Value = 1.05645161697093
$ws.Range("K7").Value = 1.057822073254539
$ws.Range("L7").Value = 1.051290483115078
$ws.Range("M7").Value = 1.067083400154881
$ws.Range("N7").Value = 1.057951899443723
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.050479134484721
$ws.Range("D8").Value = 1.054182144032022
$ws.Range("E8").Value = 1.047536940392398
$ws.Range("F8").Value = 1.063159661141581
$ws.Range("I8").Value = 1.044772310352223
$ws.Range("J8").Value = 1.055395472034386
$ws.Range("K8").Value = 1.056863802032989
$ws.Range("L8").Value = 1.050236603747463
$ws.Range("M8").Value = 1.065817380717145
$ws.Range("N8").Value = 1.056894254660228
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.047803299466382
$ws.Range("D9").Value = 1.052058810213295
$ws.Range("E9").Value = 1.045243832954212
$ws.Range("F9").Value = 1.06050055259977
$ws.Range("I9").Value = 1.043990420812896
$ws.Range("J9").Value = 1.053529654511449
$ws.Range("K9").Value = 1.055169876861557
$ws.Range("L9").Value = 1.048376548116826
$ws.Range("M9").Value = 1.063585224594569
$ws.Range("N9").Value = 1.055025787462396
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.046018958549808
$ws.Range("D10").Value = 1.050642902488887
$ws.Range("E10").Value = 1.043716224450084
$ws.Range("F10").Value = 1.058729928909636
$ws.Range("I10").Value = 1.043462805836039
$ws.Range("J10").Value = 1.052282546634136
$ws.Range("K10").Value = 1.054036986740008
$ws.Range("L10").Value = 1.04713447851965
$ws.Range("M10").Value = 1.062096234559303
$ws.Range("N10").Value = 1.053776908548852
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.045246195114106
$ws.Range("D11").Value = 1.050029710723079
$ws.Range("E11").Value = 1.0430550152381
$ws.Range("F11").Value = 1.057963721472583
$ws.Range("I11").Value = 1.043232832056209
$ws.Range("J11").Value = 1.051741761197212
$ws.Range("K11").Value = 1.053545571316718
$ws.Range("L11").Value = 1.046596160220205
$ws.Range("M11").Value = 1.061451267046618
$ws.Range("N11").Value = 1.053235355134578
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.04495913494782
$ws.Range("D12").Value = 1.049801929316951
$ws.Range("E12").Value = 1.042809450455919
$ws.Range("F12").Value = 1.057679189747669
$ws.Range("I12").Value = 1.043147181893294
$ws.Range("J12").Value = 1.051540771515126
$ws.Range("K12").Value = 1.053362906948649
$ws.Range("L12").Value = 1.04639612997981
$ws.Range("M12").Value = 1.061211662778463
$ws.Range("N12").Value = 1.05303408002409
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.045020711322575
$ws.Range("D13").Value = 1.049850789874272
$ws.Range("E13").Value = 1.042862123236311
$ws.Range("F13").Value = 1.057740219532874
$ws.Range("I13").Value = 1.043165564451735
$ws.Range("J13").Value = 1.051583889831647
$ws.Range("K13").Value = 1.053402095004493
$ws.Range("L13").Value = 1.046439040543083
$ws.Range("M13").Value = 1.061263060286878
$ws.Range("N13").Value = 1.053077259573566
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.045222467077291
$ws.Range("D14").Value = 1.050010882543411
$ws.Range("E14").Value = 1.043034716009816
$ws.Range("F14").Value = 1.057940200531449
$ws.Range("I14").Value = 1.043225756833259
$ws.Range("J14").Value = 1.05172514973688
$ws.Range("K14").Value = 1.053530474894894
$ws.Range("L14").Value = 1.046579627188877
$ws.Range("M14").Value = 1.06143146199394
$ws.Range("N14").Value = 1.053218720084067
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.045346772571059
$ws.Range("D15").Value = 1.050109518960029
$ws.Range("E15").Value = 1.04314106109926
$ws.Range("F15").Value = 1.05806342491561
$ws.Range("I15").Value = 1.043262813168355
$ws.Range("J15").Value = 1.051812168981393
$ws.Range("K15").Value = 1.053609556631856
$ws.Range("L15").Value = 1.046666237323011
$ws.Range("M15").Value = 1.061535215232247
$ws.Range("N15").Value = 1.053305862905888
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.046070241397038
$ws.Range("D16").Value = 1.050683595975583
$ws.Range("E16").Value = 1.043760112056268
$ws.Range("F16").Value = 1.058780789698276
$ws.Range("I16").Value = 1.043478036531202
$ws.Range("J16").Value = 1.052318420292735
$ws.Range("K16").Value = 1.054069582054134
$ws.Range("L16").Value = 1.047170194470143
$ws.Range("M16").Value = 1.062139034119075
$ws.Range("N16").Value = 1.05381283315216
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.046524017810212
$ws.Range("D17").Value = 1.051043673698888
$ws.Range("E17").Value = 1.04414849430557
$ws.Range("F17").Value = 1.059230902681007
$ws.Range("I17").Value = 1.043612635137969
$ws.Range("J17").Value = 1.052635769120595
$ws.Range("K17").Value = 1.05435791142152
$ws.Range("L17").Value = 1.047486180772198
$ws.Range("M17").Value = 1.062517732941438
$ws.Range("N17").Value = 1.054130632651756
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.04678868528044
$ws.Range("D18").Value = 1.051253691870639
$ws.Range("E18").Value = 1.044375055778987
$ws.Range("F18").Value = 1.059493492759299
$ws.Range("I18").Value = 1.043690998267568
$ws.Range("J18").Value = 1.052820798165153
$ws.Range("K18").Value = 1.054526005497529
$ws.Range("L18").Value = 1.047670442618578
$ws.Range("M18").Value = 1.062738599902616
$ws.Range("N18").Value = 1.05431592445878
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.04687892795326
$ws.Range("D19").Value = 1.051325301116305
$ws.Range("E19").Value = 1.044452311590762
$ws.Range("F19").Value = 1.059583037154824
$ws.Range("I19").Value = 1.043717693322852
$ws.Range("J19").Value = 1.052883875563503
$ws.Range("K19").Value = 1.05458330707796
$ws.Range("L19").Value = 1.0477332630892
$ws.Range("M19").Value = 1.062813906180344
$ws.Range("N19").Value = 1.054379091434272
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.046475333213591
$ws.Range("D20").Value = 1.051005041708784
$ws.Range("E20").Value = 1.044106822017446
$ws.Range("F20").Value = 1.05918260498452
$ws.Range("I20").Value = 1.043598209081109
$ws.Range("J20").Value = 1.052601728374122
$ws.Range("K20").Value = 1.054326985057025
$ws.Range("L20").Value = 1.047452283362252
$ws.Range("M20").Value = 1.062477104387332
$ws.Range("N20").Value = 1.054096543563519
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.045163055640118
$ws.Range("D21").Value = 1.049963739651893
$ws.Range("E21").Value = 1.042983890702292
$ws.Range("F21").Value = 1.057881309136292
$ws.Range("I21").Value = 1.043208037964521
$ws.Range("J21").Value = 1.051683555476285
$ws.Range("K21").Value = 1.053492673832597
$ws.Range("L21").Value = 1.046538229997257
$ws.Range("M21").Value = 1.061381872850632
$ws.Range("N21").Value = 1.053177066754851
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.044337849939723
$ws.Range("D22").Value = 1.049308946524258
$ws.Range("E22").Value = 1.04227807738001
$ws.Range("F22").Value = 1.057063547929421
$ws.Range("I22").Value = 1.042961404546622
$ws.Range("J22").Value = 1.051105581257193
$ws.Range("K22").Value = 1.052967352181702
$ws.Range("L22").Value = 1.045963094959387
$ws.Range("M22").Value = 1.060693056629394
$ws.Range("N22").Value = 1.052598271746074
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.044775319443496
$ws.Range("D23").Value = 1.049656072994279
$ws.Range("E23").Value = 1.04265222189935
$ws.Range("F23").Value = 1.057497019577515
$ws.Range("I23").Value = 1.043092274573724
$ws.Range("J23").Value = 1.0514120412785
$ws.Range("K23").Value = 1.053245907095048
$ws.Range("L23").Value = 1.04626802622183
$ws.Range("M23").Value = 1.06105823043529
$ws.Range("N23").Value = 1.052905166975762
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.046497331725786
$ws.Range("D24").Value = 1.051022497870021
$ws.Range("E24").Value = 1.044125651854185
$ws.Range("F24").Value = 1.059204428489714
$ws.Range("I24").Value = 1.043604728046096
$ws.Range("J24").Value = 1.052617110154134
$ws.Range("K24").Value = 1.054340959607899
$ws.Range("L24").Value = 1.047467600289595
$ws.Range("M24").Value = 1.062495462748714
$ws.Range("N24").Value = 1.054111947187424
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.048495141039477
$ws.Range("D25").Value = 1.05260780367383
$ws.Range("E25").Value = 1.045836455423188
$ws.Range("F25").Value = 1.061187620201475
$ws.Range("I25").Value = 1.044193677228202
$ws.Range("J25").Value = 1.054012579914203
$ws.Range("K25").Value = 1.055608431070803
$ws.Range("L25").Value = 1.048857773474865
$ws.Range("M25").Value = 1.0641624441667
$ws.Range("N25").Value = 1.055509398674614